$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "damit/weil/denn/um …" column to "damit/weil/denn …" ---
$ws.Range("G3").Value = "damit/weil/denn …"

# --- Column E ("Als…"): "Benutzer" / "Teilnehmer" are unified into "Nutzer" ---
$ws.Range("E4").Value  = "Nutzer"
$ws.Range("E5").Value  = "Nutzer"
$ws.Range("E6").Value  = "Nutzer"
$ws.Range("E12").Value = "Nutzer"
$ws.Range("E13").Value = "Nutzer"
$ws.Range("E14").Value = "Nutzer"
$ws.Range("E15").Value = "Nutzer"
$ws.Range("E16").Value = "Nutzer"
$ws.Range("E17").Value = "Nutzer"

# --- Column G ("möchte ich …"): reworded user-story reasons ---
$ws.Range("G4").Value  = "ich entscheiden kann, ob die Veranstaltung meinen Erwartungen entspricht"
$ws.Range("G5").Value  = "ich personalisierte funktionen nutzen kann"
$ws.Range("G6").Value  = "ich meine persoenlichen Informationen aendern kann "
$ws.Range("G7").Value  = "ich sie ueber wichtige Aenderungen informieren kann "
$ws.Range("G8").Value  = "ich die Ressourcen besser planen kann"
$ws.Range("G9").Value  = "ich die teilnehmende Personen ueber den Standort informieren kann"
$ws.Range("G10").Value = "ich alle relevanten Informationen erfassen und organisieren kann"
$ws.Range("G12").Value = "ich den Veranstaltungstermin in meinem Kalender sehen möchte "
$ws.Range("G13").Value = "ich Veranstaltungen finden möchte, die meinen Interessen entsprechen"
$ws.Range("G14").Value = "ich interessante Veranstaltungen finden möchte"
$ws.Range("G15").Value = "ich an Veranstaltungen teilnehmen kann"
$ws.Range("G16").Value = "ich planänderungen an den Veranstalter mitteilen kann"
$ws.Range("G17").Value = "sich die Qualität zukünftiger Events verbessern kann"

# --- Update the active selection to match the saved view state ---
$ws.Range("G12").Select()
